$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing lat/lon columns (D,E) out to (E,F), and the Status
# column (C) out to (D), right-to-left so nothing gets clobbered before
# it is read. Then insert the new "Date" column in C.

# Header row first, so new shared strings ("Date" then the month labels)
# land in the same order as produced by the original authoring tool.
$ws.Range("F1").Value = "lon"
$ws.Range("E1").Value = "lat"
$ws.Range("D1").Value = "Status"
$ws.Range("C1").Value = "Date"

# Row 2 - Washington DC / Modelling the World's Systems 2019
$ws.Range("F2").Value = -77.03637
$ws.Range("E2").Value = 38.89511
$ws.Range("D2").Value = "Attended"
$ws.Range("C2").Value = "May-19"

# Row 3 - Bristol / MRF AMR Annual Conference 2019
$ws.Range("F3").Value = -2.58791
$ws.Range("E3").Value = 51.454514
$ws.Range("D3").Value = "Attended"
$ws.Range("C3").Value = "Aug-19"

# Row 4 - Charleston / EPIDEMICS 2019 (status updated: now Attended)
$ws.Range("F4").Value = -79.940918
$ws.Range("E4").Value = 32.784618
$ws.Range("D4").Value = "Attended"
$ws.Range("C4").Value = "Dec-19"

# Row 5 - Edinburgh / Microbiology Society Annual Conference 2020
$ws.Range("F5").Value = -3.188267
$ws.Range("E5").Value = 55.953251
$ws.Range("D5").Value = "Will attend"
$ws.Range("C5").Value = "Mar-20"

# Row 6 - Paris / ECCMID 2020
$ws.Range("F6").Value = 2.349014
$ws.Range("E6").Value = 48.864716
$ws.Range("D6").Value = "Will attend"
$ws.Range("C6").Value = "Apr-20"

# Apply text number format to the new Date column's data cells
$ws.Range("C2:C6").NumberFormat = "@"
